# Update "paises" (countries) workbook:
#  - bump the "Datos actualizados" timestamp
#  - refresh COVID case statistics for a number of countries
#  - some countries changed rank (their row position in the table), which
#    shows up as the country label in a given row changing together with
#    the numbers for that row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / timestamp
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 23:05"

# Estados Unidos (row 4) - values refreshed, same rank
$ws.Range("B4").Value = 1851186
$ws.Range("C4").Value = 14016
$ws.Range("D4").Value = 607921
$ws.Range("E4").Value = 1136674
$ws.Range("G4").Value = 396
$ws.Range("H4").Value = 106591

# Alemania (row 12) - values refreshed, same rank
$ws.Range("B12").Value = 183744
$ws.Range("C12").Value = 250
$ws.Range("E12").Value = 9226
$ws.Range("G12").Value = 13
$ws.Range("H12").Value = 8618

# Peru overtakes Turquia -> row 13 becomes Peru (fresh numbers),
# row 14 becomes Turquia (its former, unchanged numbers)
$ws.Range("A13").Value = "Peru"
$ws.Range("B13").Value = 170039
$ws.Range("C13").Value = 5563
$ws.Range("D13").Value = 68507
$ws.Range("E13").Value = 96898
$ws.Range("G13").Value = 128
$ws.Range("H13").Value = 4634

$ws.Range("A14").Value = "Turquia"
$ws.Range("B14").Value = 164769
$ws.Range("C14").Value = 827
$ws.Range("D14").Value = 128947
$ws.Range("E14").Value = 31259
$ws.Range("G14").Value = 23
$ws.Range("H14").Value = 4563

# Camerun (row 69) - values refreshed, same rank
$ws.Range("B69").Value = 6397
$ws.Range("C69").Value = 493
$ws.Range("D69").Value = 3629
$ws.Range("E69").Value = 2569
$ws.Range("G69").Value = 8
$ws.Range("H69").Value = 199

# Sudan overtakes Guatemala -> row 72 becomes Sudan (fresh numbers),
# row 73 becomes Guatemala (its former, unchanged numbers)
$ws.Range("A72").Value = "Sudan"
$ws.Range("B72").Value = 5173
$ws.Range("C72").Value = 147
$ws.Range("D72").Value = 1522
$ws.Range("E72").Value = 3353
$ws.Range("G72").Value = 12
$ws.Range("H72").Value = 298

$ws.Range("A73").Value = "Guatemala"
$ws.Range("B73").Value = 5087
$ws.Range("C73").Value = 348
$ws.Range("D73").Value = 735
$ws.Range("E73").Value = 4244
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 108

# Uzbekistan (row 79) - values refreshed, same rank
$ws.Range("B79").Value = 3702
$ws.Range("C79").Value = 79
$ws.Range("E79").Value = 828

# Costa de Marfil overtakes Grecia -> row 83 becomes Costa de Marfil
# (fresh numbers), row 84 becomes Grecia (its former, unchanged numbers)
$ws.Range("A83").Value = "Costa de Marfil"
$ws.Range("B83").Value = 2951
$ws.Range("C83").Value = 118
$ws.Range("D83").Value = 1467
$ws.Range("E83").Value = 1451
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 33

$ws.Range("A84").Value = "Grecia"
$ws.Range("B84").Value = 2918
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 1374
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 179

# Guinea-Bisau jumps ahead of Mali, Guinea Ecuatorial and Etiopia ->
# row 107 becomes Guinea-Bisau (fresh numbers), rows 108-110 shift down
# keeping their former, unchanged numbers
$ws.Range("A107").Value = "Guinea-Bisau"
$ws.Range("B107").Value = 1339
$ws.Range("C107").Value = 83
$ws.Range("D107").Value = 53
$ws.Range("E107").Value = 1278
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 8

$ws.Range("A108").Value = "Mali"
$ws.Range("B108").Value = 1315
$ws.Range("C108").Value = 50
$ws.Range("D108").Value = 744
$ws.Range("E108").Value = 493
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 78

$ws.Range("A109").Value = "Guinea Ecuatorial"
$ws.Range("B109").Value = 1306
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 200
$ws.Range("E109").Value = 1094
$ws.Range("G109").Value = 0

$ws.Range("A110").Value = "Etiopia"
$ws.Range("B110").Value = 1257
$ws.Range("C110").Value = 85
$ws.Range("D110").Value = 217
$ws.Range("E110").Value = 1028
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 12

# Santo Tome y Principe (row 140) - values refreshed, same rank
$ws.Range("B140").Value = 484
$ws.Range("C140").Value = 1
$ws.Range("E140").Value = 404

# Ruanda (row 147) - values refreshed, same rank
$ws.Range("B147").Value = 377
$ws.Range("C147").Value = 7
$ws.Range("D147").Value = 262
$ws.Range("E147").Value = 114

# Benin (row 157) - values refreshed, same rank
$ws.Range("B157").Value = 243
$ws.Range("C157").Value = 11
$ws.Range("D157").Value = 147
$ws.Range("E157").Value = 93

# Santa Lucia overtakes Belice -> row 201 becomes Santa Lucia, row 202
# becomes Belice
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

$ws.Range("A202").Value = "Belice"
$ws.Range("D202").Value = 16
$ws.Range("H202").Value = 2
